$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

# URL (row 2, column B) - update canonical URL path segment
$meta.Range("B2").Value = "https://hl7.fr/ig/fhir/medication/StructureDefinition/fr-treatment-intent"

# Date (row 8, column B) - update generation timestamp
$meta.Range("B8").Value = "2025-05-05T08:11:38+00:00"

# --- Elements sheet updates ---
$elem = $wb.Worksheets.Item("Elements")

# Binding Value Set (row 6, column Z) - update canonical URL path segment
$elem.Range("Z6").Value = "https://hl7.fr/ig/fhir/medication/ValueSet/fr-treatment-intent"

# The column auto-fit its width slightly after the text changed (was
# 47.8125 characters, now ~48.05 characters wide) - reapply best-fit sizing.
$elem.Columns.Item(26).ColumnWidth = 47.17
